$wb = $excel.ActiveWorkbook

function Set-TextValue($range, $text) {
    # Force the cell to keep/receive a pure-text (shared string) value even
    # when the text looks like a number or a date, and restore the cell's
    # original "Normal" look (General number format, top-aligned) afterwards
    # so no visible formatting changes are introduced.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
    $range.VerticalAlignment = -4160
}

# Sheet1: MobileNumber, Date, Time, Date&Time, Enquiry_Date,
# Enquiry_PhoneNumber, User1RecId, Lead_PN, Sales_PN
$ws1 = $wb.Worksheets.Item("Sheet1")
Set-TextValue $ws1.Range("F2") "5786007259"
Set-TextValue $ws1.Range("N2") "2024-03-06"
Set-TextValue $ws1.Range("O2") "02:35:55 PM"
Set-TextValue $ws1.Range("P2") "2024-03-06 07:41:41 PM"
Set-TextValue $ws1.Range("AC2") "2024-03-06"
Set-TextValue $ws1.Range("AE2") "9963201403"
Set-TextValue $ws1.Range("AN2") "97688"
Set-TextValue $ws1.Range("AT2") "3721540577"
Set-TextValue $ws1.Range("AX2") "1870155828"

# Sheet2: MobileNumber, Enquiry_PhoneNumber, Lead_PN, Sales_PN
$ws2 = $wb.Worksheets.Item("Sheet2")
Set-TextValue $ws2.Range("F2") "5786007259"
Set-TextValue $ws2.Range("AE2") "9963201403"
Set-TextValue $ws2.Range("AT2") "3721540577"
Set-TextValue $ws2.Range("AX2") "1870155828"

# Sheet3: MobileNumber, Enquiry_PhoneNumber, Lead_PN, Sales_PN
$ws3 = $wb.Worksheets.Item("Sheet3")
Set-TextValue $ws3.Range("F2") "5786007259"
Set-TextValue $ws3.Range("AE2") "9963201403"
Set-TextValue $ws3.Range("AT2") "3721540577"
Set-TextValue $ws3.Range("AX2") "1870155828"

# Sheet4: MobileNumber, Enquiry_PhoneNumber, Lead_PN, Sales_PN
$ws4 = $wb.Worksheets.Item("Sheet4")
Set-TextValue $ws4.Range("F2") "5786007259"
Set-TextValue $ws4.Range("AE2") "9963201403"
Set-TextValue $ws4.Range("AT2") "3721540577"
Set-TextValue $ws4.Range("AX2") "1870155828"
